{"js": "// Remove the first five bullet notes from the \"Notes on Assistant Supervisor\n// Review Program\" list, leaving only the final bullet (about RV_B5 / RV_B7 /\n// DEMOGRAPHICS_ROSTER) in place. These were superseded notes that the author\n// resolved/incorporated elsewhere per the commit message.\n\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items,text\");\nawait context.sync();\n\n// Unique snippets identifying each paragraph that should be removed.\nconst snippetsToRemove = [\n  \"The RV_VILLAGE_NAME field needs to be moved up in the tree\",\n  \"The CreateAndRunPFF function is not necessary\",\n  \"In PROC RV_HOUSEHOLD_NUMBER, if the household number can\",\n  \"Should RV_B1 (line number) be automatically filled in\",\n  \"The noinputs in PROC RV_IN_HOUSEHOLD\",\n];\n\nconst items = paragraphs.items;\nfor (let i = items.length - 1; i >= 0; i--) {\n  const text = items[i].text;\n  if (snippetsToRemove.some((snippet) => text.indexOf(snippet) !== -1)) {\n    items[i].delete();\n  }\n}\n\nawait context.sync();\n", "ps1": "# Remove the first five bullet notes from the \"Notes on Assistant Supervisor\n# Review Program\" list, leaving only the final bullet (about RV_B5 / RV_B7 /\n# DEMOGRAPHICS_ROSTER) in place. These were superseded notes that the author\n# resolved/incorporated elsewhere per the commit message.\n\n$d = $word.ActiveDocument\n\n$snippets = @(\n    \"*RV_VILLAGE_NAME field needs to be moved up in the tree*\",\n    \"*CreateAndRunPFF function is not necessary*\",\n    \"*In PROC RV_HOUSEHOLD_NUMBER, if the household number can*\",\n    \"*Should RV_B1 (line number) be automatically filled in*\",\n    \"*The noinputs in PROC RV_IN_HOUSEHOLD*\"\n)\n\nfor ($i = $d.Paragraphs.Count; $i -ge 1; $i--) {\n    $p = $d.Paragraphs($i)\n    $t = $p.Range.Text\n    foreach ($snippet in $snippets) {\n        if ($t -like $snippet) {\n            $p.Range.Delete()\n            break\n        }\n    }\n}\n"}
